$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking
# strings like "30.244.96" or "1.000" are not auto-converted to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.244.96"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").Value = "1.865.41"
$ws.Range("E3").Value = "  -2.14%  "
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").Value = "234.90"
$ws.Range("E5").Value = "  -1.93%  "
$ws.Range("D6").Value = "0.9993"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").Value = "0.4671"
$ws.Range("E7").Value = "  -1.53%  "
$ws.Range("D8").Value = "0.2836"
$ws.Range("E8").Value = "  -1.24%  "
$ws.Range("D9").Value = "0.06545"
$ws.Range("E9").Value = "  -1.84%  "
$ws.Range("D10").Value = "20.85"
$ws.Range("E10").Value = "  +4.97%  "
$ws.Range("D11").Value = "0.07869"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").Value = "96.56"
$ws.Range("E12").Value = "  -4.50%  "
$ws.Range("D13").Value = "1.872.18"
$ws.Range("E13").Value = "  -1.76%  "
$ws.Range("D14").Value = "5.134"
$ws.Range("E14").Value = "  -1.10%  "
$ws.Range("D15").Value = "0.6737"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("D16").Value = "280.03"
$ws.Range("E16").Value = "  -2.14%  "
$ws.Range("D17").Value = "30.237.08"
$ws.Range("E17").Value = "  -1.35%  "
$ws.Range("D18").Value = "1.000"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "5.473"
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("D20").Value = "12.63"
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("D21").Value = "2.113.14"
$ws.Range("E21").Value = "  -2.12%  "
$ws.Range("D22").Value = "0.000007252"
$ws.Range("E22").Value = "  -3.37%  "
$ws.Range("D23").Value = "0.9989"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Value = "6.176"
$ws.Range("E24").Value = "  -1.84%  "
$ws.Range("D25").Value = "9.301"
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("D26").Value = "164.83"
$ws.Range("E26").Value = "  -1.67%  "
$ws.Range("D27").Value = "19.18"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").Value = "1.904"
$ws.Range("E28").Value = "  -6.67%  "
$ws.Range("D29").Value = "1.352"
$ws.Range("E29").Value = "  -1.95%  "
$ws.Range("D30").Value = "0.09634"
$ws.Range("E30").Value = "  -3.30%  "
$ws.Range("D31").Value = "4.400"
$ws.Range("E31").Value = "  -2.59%  "
$ws.Range("D32").Value = "1.470"
$ws.Range("E32").Value = "  -2.96%  "
$ws.Range("D33").Value = "4.102"
$ws.Range("E33").Value = "  -3.84%  "
$ws.Range("D34").Value = "0.04696"
$ws.Range("E34").Value = "  -1.56%  "
$ws.Range("D35").Value = "0.7015"
$ws.Range("E35").Value = "  -3.25%  "
$ws.Range("D36").Value = "1.100"
$ws.Range("E36").Value = "  -1.08%  "
$ws.Range("D37").Value = "2.713"
$ws.Range("E37").Value = "  -0.36%  "
$ws.Range("D38").Value = "0.01856"
$ws.Range("E38").Value = "  -2.52%  "
$ws.Range("D39").Value = "6.398"
$ws.Range("E39").Value = "  -5.37%  "
$ws.Range("D40").Value = "2.527"
$ws.Range("E40").Value = "  -3.19%  "
$ws.Range("D41").Value = "73.39"
$ws.Range("E41").Value = "  -1.28%  "
$ws.Range("D42").Value = "1.936"
$ws.Range("E42").Value = "  -2.57%  "
$ws.Range("D43").Value = "0.8455"
$ws.Range("E43").Value = "  -3.16%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.4172"
$ws.Range("E44").Value = "  -2.60%  "
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "103.98"
$ws.Range("E45").Value = "  -0.83%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "0.9993"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "7.179"
$ws.Range("E47").Value = "  -3.20%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "9.227"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "925.47"
$ws.Range("E49").Value = "  -7.10%  "
$ws.Range("D50").Value = "34.04"
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("D51").Value = "0.1132"
$ws.Range("E51").Value = "  -4.55%  "
